$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.715.83'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '3.850.04'
$ws.Range('E3').Value = '  -2.09%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'517.79"
$ws.Range('E5').Value = '  +4.90%  '
$ws.Range('D6').Value = "'140.85"
$ws.Range('E6').Value = '  -4.49%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').Value = "'0.711"
$ws.Range('E9').Value = '  -2.74%  '
$ws.Range('D10').Value = "'0.168"
$ws.Range('E10').Value = '  -4.49%  '
$ws.Range('D11').Value = "'0.0000322"
$ws.Range('E11').Value = '  -7.85%  '
$ws.Range('E12').Value = '  -4.08%  '
$ws.Range('D13').Value = "'10.31"
$ws.Range('E13').Value = '  -1.24%  '
$ws.Range('D14').Value = '4.460.75'
$ws.Range('E14').Value = '  -4.27%  '
$ws.Range('D15').Value = "'21.18"
$ws.Range('E15').Value = '  +6.59%  '
$ws.Range('D16').Value = '3.836.25'
$ws.Range('E16').Value = '  -1.88%  '
$ws.Range('D17').Value = "'13.98"
$ws.Range('E17').Value = '  -1.90%  '
$ws.Range('E18').Value = '  -2.07%  '
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('D20').Value = '68.599.40'
$ws.Range('E20').Value = '  -0.89%  '
$ws.Range('D21').Value = "'414.22"
$ws.Range('E21').Value = '  -5.66%  '
$ws.Range('D22').Value = "'3.46"
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('B23').Value = 'RenderToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D23').Value = "'12.27"
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = "'14.00"
$ws.Range('E24').Value = '  -3.59%  '
$ws.Range('D25').Value = "'86.64"
$ws.Range('E25').Value = '  -2.23%  '
$ws.Range('D26').Value = "'3.99"
$ws.Range('E26').Value = '  +4.93%  '
$ws.Range('E27').Value = '  -6.81%  '
$ws.Range('D28').Value = "'35.35"
$ws.Range('E28').Value = '  -4.72%  '
$ws.Range('D29').Value = "'13.38"
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').Value = "'677.49"
$ws.Range('E30').Value = '  -3.73%  '
$ws.Range('D31').Value = "'6.99"
$ws.Range('E31').Value = '  +14.18%  '
$ws.Range('E32').Value = '  -1.65%  '
$ws.Range('E33').Value = '  -5.28%  '
$ws.Range('D34').Value = "'66.37"
$ws.Range('E34').Value = '  +7.61%  '
$ws.Range('D35').Value = "'0.445"
$ws.Range('E35').Value = '  -5.69%  '
$ws.Range('D36').Value = '0.0₃0846'
$ws.Range('E36').Value = '  -5.98%  '
$ws.Range('D37').Value = "'39.46"
$ws.Range('E37').Value = '  -2.99%  '
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').Value = "'1.00"
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('E39').Value = '  -2.49%  '
$ws.Range('B40').Value = 'ThetaToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D40').Value = "'3.37"
$ws.Range('E40').Value = '  +9.76%  '
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('D42').Value = "'2.86"
$ws.Range('E42').Value = '  -1.75%  '
$ws.Range('E43').Value = '  -3.31%  '
$ws.Range('E44').Value = '  +4.58%  '
$ws.Range('D45').Value = "'3.44"
$ws.Range('E45').Value = '  +2.12%  '
$ws.Range('E46').Value = '  -1.47%  '
$ws.Range('D47').Value = "'0.000280"
$ws.Range('E47').Value = '  +16.84%  '
$ws.Range('E48').Value = '  -0.13%  '
$ws.Range('D49').Value = "'3.29"
$ws.Range('E49').Value = '  -2.83%  '
$ws.Range('D50').Value = "'142.90"
$ws.Range('E50').Value = '  -1.00%  '
$ws.Range('D51').Value = "'8.74"
$ws.Range('E51').Value = '  +3.13%  '
